$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1128054777"
$ws.Range("D16").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E16").Value = "1802"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 781242
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1128054777"
$ws.Range("D17").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E17").Value = "1803"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 781242
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1128054777"
$ws.Range("D18").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E18").Value = "1804"
$ws.Range("F18").Value = 27578
$ws.Range("G18").Value = 781242
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1128054777"
$ws.Range("D19").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E19").Value = "1805"
$ws.Range("F19").Value = 27578
$ws.Range("G19").Value = 781242
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1128054777"
$ws.Range("D20").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E20").Value = "1806"
$ws.Range("F20").Value = 27578
$ws.Range("G20").Value = 781242
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1128054777"
$ws.Range("D21").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E21").Value = "1807"
$ws.Range("F21").Value = 27578
$ws.Range("G21").Value = 781242
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1128054777"
$ws.Range("D22").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E22").Value = "1808"
$ws.Range("F22").Value = 27578
$ws.Range("G22").Value = 781242
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1128054777"
$ws.Range("D23").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E23").Value = "1809"
$ws.Range("F23").Value = 31249
$ws.Range("G23").Value = 781242
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1128054777"
$ws.Range("D24").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E24").Value = "1810"
$ws.Range("F24").Value = 31249
$ws.Range("G24").Value = 781242
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "45757151"
$ws.Range("D25").Value = "CLAUDIA PATRICIA PATERNINA ACEVEDO"
$ws.Range("E25").Value = "1810"
$ws.Range("F25").Value = 36000
$ws.Range("G25").Value = 900000
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1128054777"
$ws.Range("D26").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E26").Value = "1811"
$ws.Range("F26").Value = 31249
$ws.Range("G26").Value = 781242
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "45757151"
$ws.Range("D27").Value = "CLAUDIA PATRICIA PATERNINA ACEVEDO"
$ws.Range("E27").Value = "1811"
$ws.Range("F27").Value = 36000
$ws.Range("G27").Value = 900000
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1128054777"
$ws.Range("D28").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E28").Value = "1812"
$ws.Range("F28").Value = 31249
$ws.Range("G28").Value = 781242
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "45757151"
$ws.Range("D29").Value = "CLAUDIA PATRICIA PATERNINA ACEVEDO"
$ws.Range("E29").Value = "1812"
$ws.Range("F29").Value = 36000
$ws.Range("G29").Value = 900000
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1128054777"
$ws.Range("D30").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E30").Value = "1901"
$ws.Range("F30").Value = 31249
$ws.Range("G30").Value = 781242
$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "45757151"
$ws.Range("D31").Value = "CLAUDIA PATRICIA PATERNINA ACEVEDO"
$ws.Range("E31").Value = "1901"
$ws.Range("F31").Value = 36000
$ws.Range("G31").Value = 900000
$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "1128054777"
$ws.Range("D32").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E32").Value = "1902"
$ws.Range("F32").Value = 31249
$ws.Range("G32").Value = 781242
$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "45757151"
$ws.Range("D33").Value = "CLAUDIA PATRICIA PATERNINA ACEVEDO"
$ws.Range("E33").Value = "1902"
$ws.Range("F33").Value = 36000
$ws.Range("G33").Value = 900000
$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "1128054777"
$ws.Range("D34").Value = "LEWIS RAFAEL BOHORQUEZ FLOREZ"
$ws.Range("E34").Value = "1903"
$ws.Range("F34").Value = 26041
$ws.Range("G34").Value = 781242
$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "45757151"
$ws.Range("D35").Value = "CLAUDIA PATRICIA PATERNINA ACEVEDO"
$ws.Range("E35").Value = "1903"
$ws.Range("F35").Value = 30000
$ws.Range("G35").Value = 900000
